$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell used only as a formatting donor for PasteSpecial (text-percent fix-ups below).
# It is copied from and never itself modified.
$fmtDonor = $ws.Range("H2")

$ws.Range("E2").Value = "2026-03-01 05:48:15"
$ws.Range("N2").Value = "-2.1 °C 5:27 TU"
$ws.Range("O2").Value = "-0.9 °C"
$ws.Range("E3").Value = "2026-03-01 05:48:17"
$ws.Range("L3").Value = "18.7 km/h - 117º 5:13 TU"
$ws.Range("N3").Value = "-4.2 °C 5:00 TU"
$ws.Range("E4").Value = "2026-03-01 05:48:20"
$ws.Range("E5").Value = "2026-03-01 05:48:22"
$ws.Range("N5").Value = "-4.7 °C 5:27 TU"
$ws.Range("E6").Value = "2026-03-01 05:48:25"
# H6: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "89%"
$fmtDonor.Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("E7").Value = "2026-03-01 05:48:27"
$ws.Range("N7").Value = "12.8 °C 5:28 TU"
$ws.Range("O7").Value = "13.2 °C"
$ws.Range("E8").Value = "2026-03-01 05:48:29"
$ws.Range("N8").Value = "9.1 °C 5:19 TU"
$ws.Range("O8").Value = "9.4 °C"
$ws.Range("E9").Value = "2026-03-01 05:48:32"
$ws.Range("E10").Value = "2026-03-01 05:48:34"
$ws.Range("N10").Value = "4.1 °C 5:15 TU"
$ws.Range("O10").Value = "6.3 °C"
$ws.Range("E11").Value = "2026-03-01 05:48:36"
$ws.Range("E12").Value = "2026-03-01 05:48:39"
# H12: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "73%"
$fmtDonor.Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("N12").Value = "8.5 °C 5:24 TU"
$ws.Range("O12").Value = "10.4 °C"
$ws.Range("E13").Value = "2026-03-01 05:48:41"
$ws.Range("N13").Value = "3.6 °C 5:29 TU"
$ws.Range("O13").Value = "4.3 °C"
$ws.Range("E14").Value = "2026-03-01 05:48:43"
$ws.Range("N14").Value = "9.3 °C 5:29 TU"
$ws.Range("O14").Value = "11.0 °C"
$ws.Range("E15").Value = "2026-03-01 05:48:46"
$ws.Range("O15").Value = "8.5 °C"
$ws.Range("E16").Value = "2026-03-01 05:48:48"
# H16: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "83%"
$fmtDonor.Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("K16").Value = "-0.1 MJ/m2"
$ws.Range("N16").Value = "-6.6 °C 5:19 TU"
$ws.Range("O16").Value = "-5.0 °C"
$ws.Range("E17").Value = "2026-03-01 05:48:50"
$ws.Range("G17").Value = "2 cm"
$ws.Range("L17").Value = "15.5 km/h - 232º 5:05 TU"
$ws.Range("E18").Value = "2026-03-01 05:48:52"
$ws.Range("N18").Value = "5.9 °C 5:00 TU"
$ws.Range("O18").Value = "6.9 °C"
$ws.Range("E19").Value = "2026-03-01 05:48:55"
$ws.Range("N19").Value = "5.9 °C 5:00 TU"
$ws.Range("E20").Value = "2026-03-01 05:48:57"
$ws.Range("N20").Value = "-3.8 °C 5:29 TU"
$ws.Range("E21").Value = "2026-03-01 05:48:59"
$ws.Range("N21").Value = "5.7 °C 5:29 TU"
$ws.Range("O21").Value = "6.5 °C"
$ws.Range("E22").Value = "2026-03-01 05:49:02"
$ws.Range("L22").Value = "14.4 km/h - 304º 5:20 TU"
$ws.Range("N22").Value = "-6.2 °C 5:29 TU"
$ws.Range("O22").Value = "-5.3 °C"
$ws.Range("E23").Value = "2026-03-01 05:49:04"
$ws.Range("E24").Value = "2026-03-01 05:49:06"
$ws.Range("O24").Value = "4.5 °C"
$ws.Range("E25").Value = "2026-03-01 05:49:09"
# H25: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "95%"
$fmtDonor.Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("O25").Value = "-2.2 °C"
$ws.Range("E26").Value = "2026-03-01 05:49:11"
$ws.Range("J26").Value = "1025.9 hPa"
$ws.Range("N26").Value = "2.3 °C 5:19 TU"
$ws.Range("O26").Value = "2.6 °C"
$ws.Range("E27").Value = "2026-03-01 05:49:13"
$ws.Range("N27").Value = "-2.2 °C 5:19 TU"
$ws.Range("O27").Value = "-1.5 °C"
$ws.Range("E28").Value = "2026-03-01 05:49:16"
$ws.Range("N28").Value = "8.3 °C 5:10 TU"
$ws.Range("E29").Value = "2026-03-01 05:49:18"
$ws.Range("E30").Value = "2026-03-01 05:49:20"
$ws.Range("J30").Value = "1025.5 hPa"
$ws.Range("E31").Value = "2026-03-01 05:49:23"
# H31: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "63%"
$fmtDonor.Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("E32").Value = "2026-03-01 05:49:25"
$ws.Range("G32").Value = "2 cm"
$ws.Range("L32").Value = "4.7 km/h - 103º 5:11 TU"
$ws.Range("M32").Value = "4.6 °C 5:14 TU"
$ws.Range("O32").Value = "2.5 °C"
$ws.Range("E33").Value = "2026-03-01 05:49:27"
$ws.Range("E34").Value = "2026-03-01 05:49:30"
$ws.Range("N34").Value = "-0.6 °C 5:20 TU"
$ws.Range("O34").Value = "-0.2 °C"
$ws.Range("E35").Value = "2026-03-01 05:49:32"
$ws.Range("E36").Value = "2026-03-01 05:49:34"
# H36: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "76%"
$fmtDonor.Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("J36").Value = "1025.4 hPa"
$ws.Range("L36").Value = "15.1 km/h - 318º 5:24 TU"
$ws.Range("M36").Value = "12.6 °C 5:25 TU"
$ws.Range("O36").Value = "9.9 °C"
$ws.Range("E37").Value = "2026-03-01 05:49:37"
$ws.Range("N37").Value = "6.1 °C 5:05 TU"
$ws.Range("E38").Value = "2026-03-01 05:49:39"
$ws.Range("L38").Value = "6.1 km/h - 56º 5:16 TU"
$ws.Range("E39").Value = "2026-03-01 05:49:41"
$ws.Range("E40").Value = "2026-03-01 05:49:44"
# H40: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "88%"
$fmtDonor.Copy()
$ws.Range("H40").PasteSpecial(-4122)
$ws.Range("J40").Value = "1025.7 hPa"
$ws.Range("N40").Value = "5.6 °C 5:27 TU"
$ws.Range("O40").Value = "7.0 °C"
$ws.Range("E41").Value = "2026-03-01 05:49:46"
$ws.Range("J41").Value = "1025.4 hPa"
$ws.Range("E42").Value = "2026-03-01 05:49:49"
# H42: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "87%"
$fmtDonor.Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("O42").Value = "8.5 °C"
$ws.Range("E43").Value = "2026-03-01 05:49:51"
# H43: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "99%"
$fmtDonor.Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("L43").Value = "8.3 km/h - 262º 5:29 TU"
$ws.Range("N43").Value = "8.2 °C 5:00 TU"
$ws.Range("O43").Value = "8.6 °C"
$ws.Range("E44").Value = "2026-03-01 05:49:53"
# H44: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "95%"
$fmtDonor.Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("N44").Value = "-3.7 °C 5:29 TU"
$ws.Range("O44").Value = "-2.7 °C"
$ws.Range("E45").Value = "2026-03-01 05:49:55"
# H45: percent-looking text -> force text storage, then restore original cell style
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "95%"
$fmtDonor.Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("L45").Value = "13.7 km/h - 130º 5:22 TU"
$ws.Range("N45").Value = "2.5 °C 5:29 TU"
$ws.Range("O45").Value = "3.5 °C"
$ws.Range("E46").Value = "2026-03-01 05:49:58"
$ws.Range("M46").Value = "8.4 °C 5:29 TU"
